$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.898.06'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '1.882.23'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.019'
$ws.Range("E4").Value = '  +1.62%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.79'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("E6").Value = '  +1.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4688'
$ws.Range("E7").Value = '  -0.78%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3920'
$ws.Range("E8").Value = '  -1.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.04'
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07950'
$ws.Range("E10").Value = '  -1.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.013'
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.66'
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = '1.876.93'
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.951'
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.117'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.021'
$ws.Range("E16").Value = '  +1.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06794'
$ws.Range("E17").Value = '  +2.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.46'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.97'
$ws.Range("E20").Value = '  -1.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.016'
$ws.Range("E21").Value = '  +1.29%  '
$ws.Range("D22").Value = '27.913.45'
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.468'
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.361'
$ws.Range("E25").Value = '  +2.70%  '
$ws.Range("D26").Value = '2.122.11'
$ws.Range("E26").Value = '  +0.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.59'
$ws.Range("E27").Value = '  +1.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.99'
$ws.Range("E28").Value = '  -1.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.080'
$ws.Range("E29").Value = '  -1.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.464'
$ws.Range("E30").Value = '  -2.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.02'
$ws.Range("E31").Value = '  -1.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09551'
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9579'
$ws.Range("E33").Value = '  -1.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.656'
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.319'
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.347'
$ws.Range("E36").Value = '  -7.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06116'
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02238'
$ws.Range("E38").Value = '  -1.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.211'
$ws.Range("E39").Value = '  -1.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.137'
$ws.Range("E40").Value = '  -0.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5904'
$ws.Range("E41").Value = '  -1.84%  '
$ws.Range("E42").Value = '  -1.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.20'
$ws.Range("E43").Value = '  -1.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.270'
$ws.Range("E44").Value = '  +0.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5651'
$ws.Range("E45").Value = '  -1.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.14'
$ws.Range("E46").Value = '  -1.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.406'
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.918'
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06859'
$ws.Range("E49").Value = '  +0.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '113.56'
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("E51").Value = '  -1.11%  '
